$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The profile dataset was regenerated (for a better r2 score) and truncated
# from 60 samples down to 50 samples, so remove the now-unused trailing rows.
$ws.Range("A52:E61").EntireRow.Delete() | Out-Null

# Write the updated Height / Voltage / Current / Average Height values
# for each of the 50 remaining samples (rows 2-51).
$ws.Cells.Item(2, 2).Value = 2.095160230104522
$ws.Cells.Item(2, 3).Value = 5.005965860545533
$ws.Cells.Item(2, 4).Value = 18.00357951632732
$ws.Cells.Item(2, 5).Value = 2.21447744101518
$ws.Cells.Item(3, 2).Value = 1.886232637392391
$ws.Cells.Item(3, 3).Value = 5.01641224018114
$ws.Cells.Item(3, 4).Value = 18.00984734410869
$ws.Cells.Item(3, 5).Value = 2.21447744101518
$ws.Cells.Item(4, 2).Value = 2.19309545946011
$ws.Cells.Item(4, 3).Value = 5.001069099077753
$ws.Cells.Item(4, 4).Value = 18.00064145944665
$ws.Cells.Item(4, 5).Value = 2.21447744101518
$ws.Cells.Item(5, 2).Value = 2.473148346168887
$ws.Cells.Item(5, 3).Value = 4.987066454742314
$ws.Cells.Item(5, 4).Value = 17.99223987284539
$ws.Cells.Item(5, 5).Value = 2.21447744101518
$ws.Cells.Item(6, 2).Value = 1.986185746177197
$ws.Cells.Item(6, 3).Value = 5.011414584741899
$ws.Cells.Item(6, 4).Value = 18.00684875084514
$ws.Cells.Item(6, 5).Value = 2.21447744101518
$ws.Cells.Item(7, 2).Value = 1.672671793337828
$ws.Cells.Item(7, 3).Value = 5.027090282383868
$ws.Cells.Item(7, 4).Value = 18.01625416943032
$ws.Cells.Item(7, 5).Value = 2.21447744101518
$ws.Cells.Item(8, 2).Value = 1.842194708816258
$ws.Cells.Item(8, 3).Value = 5.018614136609946
$ws.Cells.Item(8, 4).Value = 18.01116848196597
$ws.Cells.Item(8, 5).Value = 2.21447744101518
$ws.Cells.Item(9, 2).Value = 1.866815088720299
$ws.Cells.Item(9, 3).Value = 5.017383117614744
$ws.Cells.Item(9, 4).Value = 18.01042987056885
$ws.Cells.Item(9, 5).Value = 2.21447744101518
$ws.Cells.Item(10, 2).Value = 1.989392448091735
$ws.Cells.Item(10, 3).Value = 5.011254249646172
$ws.Cells.Item(10, 4).Value = 18.0067525497877
$ws.Cells.Item(10, 5).Value = 2.21447744101518
$ws.Cells.Item(11, 2).Value = 1.856501497959063
$ws.Cells.Item(11, 3).Value = 5.017898797152806
$ws.Cells.Item(11, 4).Value = 18.01073927829168
$ws.Cells.Item(11, 5).Value = 2.21447744101518
$ws.Cells.Item(12, 2).Value = 2.269903865874106
$ws.Cells.Item(12, 3).Value = 4.997228678757054
$ws.Cells.Item(12, 4).Value = 17.99833720725423
$ws.Cells.Item(12, 5).Value = 2.21447744101518
$ws.Cells.Item(13, 2).Value = 2.447936313142641
$ws.Cells.Item(13, 3).Value = 4.988327056393627
$ws.Cells.Item(13, 4).Value = 17.99299623383618
$ws.Cells.Item(13, 5).Value = 2.21447744101518
$ws.Cells.Item(14, 2).Value = 2.349364768306158
$ws.Cells.Item(14, 3).Value = 4.993255633635451
$ws.Cells.Item(14, 4).Value = 17.99595338018127
$ws.Cells.Item(14, 5).Value = 2.21447744101518
$ws.Cells.Item(15, 2).Value = 2.571803735167741
$ws.Cells.Item(15, 3).Value = 4.982133685292372
$ws.Cells.Item(15, 4).Value = 17.98928021117542
$ws.Cells.Item(15, 5).Value = 2.21447744101518
$ws.Cells.Item(16, 2).Value = 2.450790398225022
$ws.Cells.Item(16, 3).Value = 4.988184352139508
$ws.Cells.Item(16, 4).Value = 17.9929106112837
$ws.Cells.Item(16, 5).Value = 2.21447744101518
$ws.Cells.Item(17, 2).Value = 2.177408003559448
$ws.Cells.Item(17, 3).Value = 5.001853471872787
$ws.Cells.Item(17, 4).Value = 18.00111208312367
$ws.Cells.Item(17, 5).Value = 2.21447744101518
$ws.Cells.Item(18, 2).Value = 2.149392733856007
$ws.Cells.Item(18, 3).Value = 5.003254235357959
$ws.Cells.Item(18, 4).Value = 18.00195254121477
$ws.Cells.Item(18, 5).Value = 2.21447744101518
$ws.Cells.Item(19, 2).Value = 1.948349470745689
$ws.Cells.Item(19, 3).Value = 5.013306398513475
$ws.Cells.Item(19, 4).Value = 18.00798383910809
$ws.Cells.Item(19, 5).Value = 2.21447744101518
$ws.Cells.Item(20, 2).Value = 2.279020965528495
$ws.Cells.Item(20, 3).Value = 4.996772823774334
$ws.Cells.Item(20, 4).Value = 17.9980636942646
$ws.Cells.Item(20, 5).Value = 2.21447744101518
$ws.Cells.Item(21, 2).Value = 2.14037559436109
$ws.Cells.Item(21, 3).Value = 5.003705092332704
$ws.Cells.Item(21, 4).Value = 18.00222305539962
$ws.Cells.Item(21, 5).Value = 2.21447744101518
$ws.Cells.Item(22, 2).Value = 2.034771198281994
$ws.Cells.Item(22, 3).Value = 5.00898531213666
$ws.Cells.Item(22, 4).Value = 18.00539118728199
$ws.Cells.Item(22, 5).Value = 2.21447744101518
$ws.Cells.Item(23, 2).Value = 1.920614179480678
$ws.Cells.Item(23, 3).Value = 5.014693163076725
$ws.Cells.Item(23, 4).Value = 18.00881589784603
$ws.Cells.Item(23, 5).Value = 2.21447744101518
$ws.Cells.Item(24, 2).Value = 2.431809048662672
$ws.Cells.Item(24, 3).Value = 4.989133419617626
$ws.Cells.Item(24, 4).Value = 17.99348005177058
$ws.Cells.Item(24, 5).Value = 2.21447744101518
$ws.Cells.Item(25, 2).Value = 2.296659729759078
$ws.Cells.Item(25, 3).Value = 4.995890885562805
$ws.Cells.Item(25, 4).Value = 17.99753453133768
$ws.Cells.Item(25, 5).Value = 2.21447744101518
$ws.Cells.Item(26, 2).Value = 2.124126023932716
$ws.Cells.Item(26, 3).Value = 5.004517570854123
$ws.Cells.Item(26, 4).Value = 18.00271054251247
$ws.Cells.Item(26, 5).Value = 2.21447744101518
$ws.Cells.Item(27, 2).Value = 2.010646543334365
$ws.Cells.Item(27, 3).Value = 5.010191544884041
$ws.Cells.Item(27, 4).Value = 18.00611492693043
$ws.Cells.Item(27, 5).Value = 2.21447744101518
$ws.Cells.Item(28, 2).Value = 2.208179516372375
$ws.Cells.Item(28, 3).Value = 5.00031489623214
$ws.Cells.Item(28, 4).Value = 18.00018893773928
$ws.Cells.Item(28, 5).Value = 2.21447744101518
$ws.Cells.Item(29, 2).Value = 2.268520235013024
$ws.Cells.Item(29, 3).Value = 4.997297860300108
$ws.Cells.Item(29, 4).Value = 17.99837871618006
$ws.Cells.Item(29, 5).Value = 2.21447744101518
$ws.Cells.Item(30, 2).Value = 2.310180394542469
$ws.Cells.Item(30, 3).Value = 4.995214852323635
$ws.Cells.Item(30, 4).Value = 17.99712891139418
$ws.Cells.Item(30, 5).Value = 2.21447744101518
$ws.Cells.Item(31, 2).Value = 2.333468508036785
$ws.Cells.Item(31, 3).Value = 4.99405044664892
$ws.Cells.Item(31, 4).Value = 17.99643026798935
$ws.Cells.Item(31, 5).Value = 2.21447744101518
$ws.Cells.Item(32, 2).Value = 2.985956633686791
$ws.Cells.Item(32, 3).Value = 4.96142604036642
$ws.Cells.Item(32, 4).Value = 17.97685562421985
$ws.Cells.Item(32, 5).Value = 2.21447744101518
$ws.Cells.Item(33, 2).Value = 2.636029985805414
$ws.Cells.Item(33, 3).Value = 4.978922372760488
$ws.Cells.Item(33, 4).Value = 17.98735342365629
$ws.Cells.Item(33, 5).Value = 2.21447744101518
$ws.Cells.Item(34, 2).Value = 2.302688747356581
$ws.Cells.Item(34, 3).Value = 4.99558943468293
$ws.Cells.Item(34, 4).Value = 17.99735366080976
$ws.Cells.Item(34, 5).Value = 2.21447744101518
$ws.Cells.Item(35, 2).Value = 2.016652872538359
$ws.Cells.Item(35, 3).Value = 5.009891228423841
$ws.Cells.Item(35, 4).Value = 18.0059347370543
$ws.Cells.Item(35, 5).Value = 2.21447744101518
$ws.Cells.Item(36, 2).Value = 2.17372797329359
$ws.Cells.Item(36, 3).Value = 5.002037473386079
$ws.Cells.Item(36, 4).Value = 18.00122248403165
$ws.Cells.Item(36, 5).Value = 2.21447744101518
$ws.Cells.Item(37, 2).Value = 2.541973279656366
$ws.Cells.Item(37, 3).Value = 4.983625208067941
$ws.Cells.Item(37, 4).Value = 17.99017512484076
$ws.Cells.Item(37, 5).Value = 2.21447744101518
$ws.Cells.Item(38, 2).Value = 2.381169862853769
$ws.Cells.Item(38, 3).Value = 4.99166537890807
$ws.Cells.Item(38, 4).Value = 17.99499922734484
$ws.Cells.Item(38, 5).Value = 2.21447744101518
$ws.Cells.Item(39, 2).Value = 2.100603218166251
$ws.Cells.Item(39, 3).Value = 5.005693711142446
$ws.Cells.Item(39, 4).Value = 18.00341622668547
$ws.Cells.Item(39, 5).Value = 2.21447744101518
$ws.Cells.Item(40, 2).Value = 1.884052263777862
$ws.Cells.Item(40, 3).Value = 5.016521258861866
$ws.Cells.Item(40, 4).Value = 18.00991275531712
$ws.Cells.Item(40, 5).Value = 2.21447744101518
$ws.Cells.Item(41, 2).Value = 2.161764804657714
$ws.Cells.Item(41, 3).Value = 5.002635631817873
$ws.Cells.Item(41, 4).Value = 18.00158137909072
$ws.Cells.Item(41, 5).Value = 2.21447744101518
$ws.Cells.Item(42, 2).Value = 2.392208627605666
$ws.Cells.Item(42, 3).Value = 4.991113440670476
$ws.Cells.Item(42, 4).Value = 17.99466806440229
$ws.Cells.Item(42, 5).Value = 2.21447744101518
$ws.Cells.Item(43, 2).Value = 2.408168283382676
$ws.Cells.Item(43, 3).Value = 4.990315457881625
$ws.Cells.Item(43, 4).Value = 17.99418927472897
$ws.Cells.Item(43, 5).Value = 2.21447744101518
$ws.Cells.Item(44, 2).Value = 2.184097045171919
$ws.Cells.Item(44, 3).Value = 5.001519019792163
$ws.Cells.Item(44, 4).Value = 18.0009114118753
$ws.Cells.Item(44, 5).Value = 2.21447744101518
$ws.Cells.Item(45, 2).Value = 2.120199862418421
$ws.Cells.Item(45, 3).Value = 5.004713878929838
$ws.Cells.Item(45, 4).Value = 18.0028283273579
$ws.Cells.Item(45, 5).Value = 2.21447744101518
$ws.Cells.Item(46, 2).Value = 2.140169325256095
$ws.Cells.Item(46, 3).Value = 5.003715405787954
$ws.Cells.Item(46, 4).Value = 18.00222924347277
$ws.Cells.Item(46, 5).Value = 2.21447744101518
$ws.Cells.Item(47, 2).Value = 2.186425287314723
$ws.Cells.Item(47, 3).Value = 5.001402607685023
$ws.Cells.Item(47, 4).Value = 18.00084156461102
$ws.Cells.Item(47, 5).Value = 2.21447744101518
$ws.Cells.Item(48, 2).Value = 2.450161665993952
$ws.Cells.Item(48, 3).Value = 4.988215788751061
$ws.Cells.Item(48, 4).Value = 17.99292947325064
$ws.Cells.Item(48, 5).Value = 2.21447744101518
$ws.Cells.Item(49, 2).Value = 2.310948709160934
$ws.Cells.Item(49, 3).Value = 4.995176436592712
$ws.Cells.Item(49, 4).Value = 17.99710586195563
$ws.Cells.Item(49, 5).Value = 2.21447744101518
$ws.Cells.Item(50, 2).Value = 2.427055534384369
$ws.Cells.Item(50, 3).Value = 4.98937109533154
$ws.Cells.Item(50, 4).Value = 17.99362265719892
$ws.Cells.Item(50, 5).Value = 2.21447744101518
$ws.Cells.Item(51, 2).Value = 2.335098885866699
$ws.Cells.Item(51, 3).Value = 4.993968927757424
$ws.Cells.Item(51, 4).Value = 17.99638135665445
$ws.Cells.Item(51, 5).Value = 2.21447744101518
